# Correcting market share tab for updated scenario 3s
$wb = $excel.ActiveWorkbook

$wsPlatform = $wb.Worksheets.Item("Platform Coverage")
$wsMarket   = $wb.Worksheets.Item("MarketShare")

# --- Data fix on MarketShare: move the 2026-2040 "1"s from row 3 (Old Product B)
# to row 2 (New Product A), leaving row 3 with only 2018-2025 populated.
$wsMarket.Range("L2:Z2").Value = 1
$wsMarket.Range("L3:Z3").ClearContents()

# --- Selection bookkeeping: MarketShare's selection moves onto the newly
# filled range, and MarketShare becomes the active (visible/selected) tab,
# matching the diff where tabSelected moves from "Platform Coverage" to
# "MarketShare".
$wsMarket.Activate()
$wsMarket.Range("L2:Z2").Select() | Out-Null
